# Applies the "Corrects mistakes (and changes of mind)" edit to the
# variable_correspondence workbook:
#   - inserts a new "socioEnvContextOther" row right after "socioEnvContext"
#   - inserts a new "methodology.analysesOther" row right after "methodology.analyses"
#   - flips the "Biological samples" row's vote from NO to YES and trims its
#     motivation text
#   - shrinks a couple of column widths to match the new content
#   - leaves the cursor on F14 with the view scrolled back to the top-left

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "socioEnvContextOther" right after "socioEnvContext" (old row 5) ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value2 = "socioEnvContextOther"
$ws.Range("D5").Value2 = "N"
$ws.Range("F5").Value2 = "YES"
$ws.Range("L5").Value2 = "Complement socioEnvContext"

# --- Insert "methodology.analysesOther" right after "methodology.analyses" ---
# (methodology.analyses is now on row 6 after the previous insert)
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value2 = "methodology.analysesOther"
$ws.Range("D7").Value2 = "AQ"
$ws.Range("F7").Value2 = "YES"
$ws.Range("L7").Value2 = "Complement methodology.analyses"

# --- "Biological samples" row moved from row 14 to row 16 by the two inserts;
#      the vote changes from NO to YES and the motivation text is shortened ---
$ws.Range("F16").Value2 = "YES"
$ws.Range("L16").Value2 = "Could be of interest although maybe not one of the most important ones"

# --- Column width tweaks (B, C, L got slightly narrower) ---
$ws.Columns.Item(2).ColumnWidth = 38.83333333333333
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Columns.Item(12).ColumnWidth = 126.83333333333333

# --- Restore view: scroll back to A1 and leave the cursor on F14 ---
$ws.Range("A1").Select()
$ws.Range("F14").Select()
